$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in columns A and B first, row by row (matches author's original entry order)
$ws.Range("A26").Value = "T6767000"
$ws.Range("B26").Value = "WHAT IS CURRENT HIGHEST GRADE EVER COMPLETED? 2011"

$ws.Range("A27").Value = "U1718000"
$ws.Range("B27").Value = "TYPE OF BUS OR INDUSTRY CODE (2002 CENSUS 4-DIGIT)"

$ws.Range("A28").Value = "R0515100"
$ws.Range("B28").Value = "PERCENT CHANCE R HAS COLLEGE DEGREE BY 30 YEARS OLD 1997"

$ws.Range("A29").Value = "U1719400"
$ws.Range("B29").Value = "OCCUPATION/JOB CODE (2002 CENSUS 4-DIGIT) 2017"

# Then fill in column C (variable names) afterward
$ws.Range("C26").Value = "high.grade.completed.11"
$ws.Range("C27").Value = "type.business"
$ws.Range("C28").Value = "pct.college.30.1997"
$ws.Range("C29").Value = "occupation"

$ws.Range("C30").Select()
